$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# 1. EN fluency question wording tweak: "reasonably fluent" -> "you could
#    effectively communicate". Rewrite the whole sentence (preserving every
#    surrounding run's own formatting/rsid) so only the target run's text
#    actually changes, instead of letting identically-formatted neighboring
#    runs coalesce into one run.
$rngFluent = $d.Content
$rngFluent.Find.Execute("How old were you when you felt reasonably fluent in English? __________")
$fluentRange = $d.Range($rngFluent.Start, $rngFluent.End)
$fluentXml = "<w:p $wNs>" + `
  "<w:r w:rsidRPr='00D15D28'><w:rPr><w:rFonts w:asciiTheme='minorBidi' w:hAnsiTheme='minorBidi' w:cstheme='minorBidi'/><w:color w:val='000000' w:themeColor='text1'/></w:rPr><w:t xml:space='preserve'>How old were you when you </w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:asciiTheme='minorBidi' w:hAnsiTheme='minorBidi' w:cstheme='minorBidi'/><w:color w:val='000000' w:themeColor='text1'/></w:rPr><w:t xml:space='preserve'>felt </w:t></w:r>" + `
  "<w:r w:rsidR='00AF1C1D'><w:rPr><w:rFonts w:asciiTheme='minorBidi' w:hAnsiTheme='minorBidi' w:cstheme='minorBidi'/><w:color w:val='000000' w:themeColor='text1'/></w:rPr><w:t>you could effectively communicate</w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:asciiTheme='minorBidi' w:hAnsiTheme='minorBidi' w:cstheme='minorBidi'/><w:color w:val='000000' w:themeColor='text1'/></w:rPr><w:t xml:space='preserve'> in</w:t></w:r>" + `
  "<w:r w:rsidRPr='00D15D28'><w:rPr><w:rFonts w:asciiTheme='minorBidi' w:hAnsiTheme='minorBidi' w:cstheme='minorBidi'/><w:color w:val='000000' w:themeColor='text1'/></w:rPr><w:t xml:space='preserve'> English? __________</w:t></w:r>" + `
  "</w:p>"
$fluentRange.InsertXML($fluentXml)

# 2. Move <w:lastRenderedPageBreak/> from the "Poor" run up to the
#    "Which social class group do you identify with?" run that precedes it
#    (repagination after the text tweak above moved where the page breaks).

# 2a. Drop it from the "Poor" run.
$rngPoor = $d.Content
$rngPoor.Find.Execute("Poor")
$poorRange = $d.Range($rngPoor.Start, $rngPoor.End)
$poorXml = "<w:p $wNs>" + `
           "<w:r w:rsidRPr='00D15D28'>" + `
           "<w:rPr><w:rFonts w:asciiTheme='minorBidi' w:eastAsia='Poppins' w:hAnsiTheme='minorBidi' w:cstheme='minorBidi'/></w:rPr>" + `
           "<w:t>Poor</w:t>" + `
           "</w:r></w:p>"
$poorRange.InsertXML($poorXml)

# 2b. Add it to the "Which social class group do you identify with?" run.
$rngClass = $d.Content
$rngClass.Find.Execute("Which social class group do you identify with?")
$classRange = $d.Range($rngClass.Start, $rngClass.End)
$classXml = "<w:p $wNs>" + `
            "<w:r w:rsidRPr='00D15D28'>" + `
            "<w:rPr><w:rFonts w:asciiTheme='minorBidi' w:eastAsia='Poppins' w:hAnsiTheme='minorBidi' w:cstheme='minorBidi'/></w:rPr>" + `
            "<w:lastRenderedPageBreak/>" + `
            "<w:t>Which social class group do you identify with?</w:t>" + `
            "</w:r></w:p>"
$classRange.InsertXML($classXml)
